$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 24,14
$arr[0,0] = 1.357478060742437
$arr[0,1] = 0.1530760410243772
$arr[0,2] = 0.1194354720364359
$arr[0,3] = 0.1245758532213168
$arr[0,4] = 1.752131314218971
$arr[0,5] = 0
$arr[0,6] = 0.07973214163530429
$arr[0,7] = 0
$arr[0,8] = 0.161299038329541
$arr[0,9] = 0
$arr[0,10] = 0.2744277187219382
$arr[0,11] = 0.3184569435170133
$arr[0,12] = 1.7218215251701
$arr[0,13] = 4.575683049327694
$arr[1,0] = 1.279850709967377
$arr[1,1] = 0.1432450312874494
$arr[1,2] = 0.1192628499359962
$arr[1,3] = 0.125268345086099
$arr[1,4] = 1.756729400771874
$arr[1,5] = 0
$arr[1,6] = 0.07973214163530429
$arr[1,7] = 0
$arr[1,8] = 0.1622816523651527
$arr[1,9] = 0
$arr[1,10] = 0.2715942972574581
$arr[1,11] = 0.3060752846758419
$arr[1,12] = 1.738201905344983
$arr[1,13] = 4.589442455110571
$arr[2,0] = 1.232583085629301
$arr[2,1] = 0.137161272822425
$arr[2,2] = 0.1191855449299304
$arr[2,3] = 0.1257211810429499
$arr[2,4] = 1.760423353963603
$arr[2,5] = 0
$arr[2,6] = 0.07973214163530429
$arr[2,7] = 0
$arr[2,8] = 0.1629171745434213
$arr[2,9] = 0
$arr[2,10] = 0.2699455850554813
$arr[2,11] = 0.2985791830747075
$arr[2,12] = 1.748849901457888
$arr[2,13] = 4.600302111763824
$arr[3,0] = 1.213422029166651
$arr[3,1] = 0.1346702756366653
$arr[3,2] = 0.1191612859694118
$arr[3,3] = 0.1259126794800829
$arr[3,4] = 1.762147729475636
$arr[3,5] = 0
$arr[3,6] = 0.07973214163530429
$arr[3,7] = 0
$arr[3,8] = 0.1631842639411922
$arr[3,9] = 0
$arr[3,10] = 0.2692967146946259
$arr[3,11] = 0.2955514316987973
$arr[3,12] = 1.753337581230131
$arr[3,13] = 4.605333945374014
$arr[4,0] = 1.210246482883036
$arr[4,1] = 0.1342559374007664
$arr[4,2] = 0.1191576961640237
$arr[4,3] = 0.1259448986475755
$arr[4,4] = 1.762447294712274
$arr[4,5] = 0
$arr[4,6] = 0.07973214163530429
$arr[4,7] = 0
$arr[4,8] = 0.1632291040971876
$arr[4,9] = 0
$arr[4,10] = 0.2691903621091924
$arr[4,11] = 0.2950503120599279
$arr[4,12] = 1.754091730599622
$arr[4,13] = 4.606206107381212
$arr[5,0] = 1.232324262673018
$arr[5,1] = 0.1371277260433459
$arr[5,2] = 0.1191851883962975
$arr[5,3] = 0.1257237354431089
$arr[5,4] = 1.76044572238294
$arr[5,5] = 0
$arr[5,6] = 0.07973214163530429
$arr[5,7] = 0
$arr[5,8] = 0.1629207437535092
$arr[5,9] = 0
$arr[5,10] = 0.269936740923626
$arr[5,11] = 0.2985382402327161
$arr[5,12] = 1.748909822357838
$arr[5,13] = 4.600367517287083
$arr[6,0] = 1.330630759737573
$arr[6,1] = 0.1496962338580232
$arr[6,2] = 0.1193700161123807
$arr[6,3] = 0.1248088949287354
$arr[6,4] = 1.753536082980425
$arr[6,5] = 0
$arr[6,6] = 0.07973214163530429
$arr[6,7] = 0
$arr[6,8] = 0.1616311702708058
$arr[6,9] = 0
$arr[6,10] = 0.2734319256642124
$arr[6,11] = 0.3141658296617535
$arr[6,12] = 1.727346996114779
$arr[6,13] = 4.579926865676526
$arr[7,0] = 1.526502347137637
$arr[7,1] = 0.1739620133617734
$arr[7,2] = 0.1199588316441051
$arr[7,3] = 0.1232336408118266
$arr[7,4] = 1.746891606093811
$arr[7,5] = 0
$arr[7,6] = 0.07973214163530429
$arr[7,7] = 0
$arr[7,8] = 0.159357114871824
$arr[7,9] = 0
$arr[7,10] = 0.2810040988719749
$arr[7,11] = 0.3456458571932828
$arr[7,12] = 1.689742833054538
$arr[7,13] = 4.558976016853336
$arr[8,0] = 1.672245089164278
$arr[8,1] = 0.191553616241265
$arr[8,2] = 0.1205278997622869
$arr[8,3] = 0.1222088324085897
$arr[8,4] = 1.746217339310931
$arr[8,5] = 0
$arr[8,6] = 0.07973214163530429
$arr[8,7] = 0
$arr[8,8] = 0.157840767748525
$arr[8,9] = 0
$arr[8,10] = 0.2870005086179646
$arr[8,11] = 0.369273430470983
$arr[8,12] = 1.664962582492912
$arr[8,13] = 4.555254460797215
$arr[9,0] = 1.738936320920345
$arr[9,1] = 0.1995043562099568
$arr[9,2] = 0.1208161000314618
$arr[9,3] = 0.1217712239271433
$arr[9,4] = 1.746823906355829
$arr[9,5] = 0
$arr[9,6] = 0.07973214163530429
$arr[9,7] = 0
$arr[9,8] = 0.1571842709706903
$arr[9,9] = 0
$arr[9,10] = 0.2898215493158602
$arr[9,11] = 0.3801287789972534
$arr[9,12] = 1.654306541597848
$arr[9,13] = 4.556097645374962
$arr[10,0] = 1.764245798480147
$arr[10,1] = 0.2025075466328019
$arr[10,2] = 0.120929422833548
$arr[10,3] = 0.1216096098854189
$arr[10,4] = 1.747184865425695
$arr[10,5] = 0
$arr[10,6] = 0.07973214163530429
$arr[10,7] = 0
$arr[10,8] = 0.1569404467191866
$arr[10,9] = 0
$arr[10,10] = 0.2909031176539401
$arr[10,11] = 0.3842545988468089
$arr[10,12] = 1.650359984906693
$arr[10,13] = 4.556781690891057
$arr[11,0] = 1.758792528258482
$arr[11,1] = 0.2018610952586641
$arr[11,2] = 0.1209048309685485
$arr[11,3] = 0.1216442342445632
$arr[11,4] = 1.74710128941139
$arr[11,5] = 0
$arr[11,6] = 0.07973214163530429
$arr[11,7] = 0
$arr[11,8] = 0.1569927463617016
$arr[11,9] = 0
$arr[11,10] = 0.2906695926747886
$arr[11,11] = 0.3833653606980292
$arr[11,12] = 1.651206003743496
$arr[11,13] = 4.556618146208137
$arr[12,0] = 1.741017453675568
$arr[12,1] = 0.1997515832424881
$arr[12,2] = 0.1208253394108709
$arr[12,3] = 0.1217578457609623
$arr[12,4] = 1.746850972297665
$arr[12,5] = 0
$arr[12,6] = 0.07973214163530429
$arr[12,7] = 0
$arr[12,8] = 0.1571641157320105
$arr[12,9] = 0
$arr[12,10] = 0.2899102646248224
$arr[12,11] = 0.3804679108300704
$arr[12,12] = 1.653980079298258
$arr[12,13] = 4.556146611760454
$arr[13,0] = 1.730136821615588
$arr[13,1] = 0.1984584546983115
$arr[13,2] = 0.1207771930069796
$arr[13,3] = 0.1218279695460973
$arr[13,4] = 1.746714738516374
$arr[13,5] = 0
$arr[13,6] = 0.07973214163530429
$arr[13,7] = 0
$arr[13,8] = 0.1572697061724453
$arr[13,9] = 0
$arr[13,10] = 0.2894468833427339
$arr[13,11] = 0.3786951024843646
$arr[13,12] = 1.655690825783871
$arr[13,13] = 4.555905286154939
$arr[14,0] = 1.667894428470504
$arr[14,1] = 0.1910329631686807
$arr[14,2] = 0.120509652811009
$arr[14,3] = 0.1222380053301944
$arr[14,4] = 1.746196070959456
$arr[14,5] = 0
$arr[14,6] = 0.07973214163530429
$arr[14,7] = 0
$arr[14,8] = 0.1578843404150918
$arr[14,9] = 0
$arr[14,10] = 0.2868180141338996
$arr[14,11] = 0.3685661398856155
$arr[14,12] = 1.665671386456545
$arr[14,13] = 4.555250395097033
$arr[15,0] = 1.629810149942443
$arr[15,1] = 0.1864643072149761
$arr[15,2] = 0.1203530160379529
$arr[15,3] = 0.1224968614991333
$arr[15,4] = 1.746111755306018
$arr[15,5] = 0
$arr[15,6] = 0.07973214163530429
$arr[15,7] = 0
$arr[15,8] = 0.1582699187598462
$arr[15,9] = 0
$arr[15,10] = 0.2852290938580353
$arr[15,11] = 0.3623795827904885
$arr[15,12] = 1.67195205657859
$arr[15,13] = 4.555498246597153
$arr[16,0] = 1.607942074556036
$arr[16,1] = 0.1838316726939695
$arr[16,2] = 0.1202656850036163
$arr[16,3] = 0.1226484396630845
$arr[16,4] = 1.746149226654012
$arr[16,5] = 0
$arr[16,6] = 0.07973214163530429
$arr[16,7] = 0
$arr[16,8] = 0.1584948279743981
$arr[16,9] = 0
$arr[16,10] = 0.2843239718655894
$arr[16,11] = 0.3588313299504549
$arr[16,12] = 1.675622571627422
$arr[16,13] = 4.555879514251529
$arr[17,0] = 1.600544317460901
$arr[17,1] = 0.1829394766156156
$arr[17,2] = 0.1202365915505297
$arr[17,3] = 0.1227002239412185
$arr[17,4] = 1.746176680802236
$arr[17,5] = 0
$arr[17,6] = 0.07973214163530429
$arr[17,7] = 0
$arr[17,8] = 0.1585715171109205
$arr[17,9] = 0
$arr[17,10] = 0.2840190248483054
$arr[17,11] = 0.3576316935522144
$arr[17,12] = 1.676875313843368
$arr[17,13] = 4.556049601909621
$arr[18,0] = 1.633860468010823
$arr[18,1] = 0.1869511528948351
$arr[18,2] = 0.1203694046186712
$arr[18,3] = 0.1224690273975675
$arr[18,4] = 1.746111834010975
$arr[18,5] = 0
$arr[18,6] = 0.07973214163530429
$arr[18,7] = 0
$arr[18,8] = 0.1582285489339244
$arr[18,9] = 0
$arr[18,10] = 0.285397328909923
$arr[18,11] = 0.3630371098282339
$arr[18,12] = 1.671277461744999
$arr[18,13] = 4.555447156376744
$arr[19,0] = 1.746236943879467
$arr[19,1] = 0.200371404935396
$arr[19,2] = 0.1208485745942411
$arr[19,3] = 0.1217243641540016
$arr[19,4] = 1.746920934519437
$arr[19,5] = 0
$arr[19,6] = 0.07973214163530429
$arr[19,7] = 0
$arr[19,8] = 0.1571136508413593
$arr[19,9] = 0
$arr[19,10] = 0.2901329375325048
$arr[19,11] = 0.3813185533892067
$arr[19,12] = 1.653162860092735
$arr[19,13] = 4.556275213204657
$arr[20,0] = 1.820001126835791
$arr[20,1] = 0.2090980564121026
$arr[20,2] = 0.121186134769026
$arr[20,3] = 0.121261569406264
$arr[20,4] = 1.748214812819157
$arr[20,5] = 0
$arr[20,6] = 0.07973214163530429
$arr[20,7] = 0
$arr[20,8] = 0.1564128383781567
$arr[20,9] = 0
$arr[20,10] = 0.2933054131346182
$arr[20,11] = 0.3933546321384469
$arr[20,12] = 1.641840664693035
$arr[20,13] = 4.558942463405515
$arr[21,0] = 1.780603017806925
$arr[21,1] = 0.2044445735256204
$arr[21,2] = 0.1210037502182857
$arr[21,3] = 0.121506389828328
$arr[21,4] = 1.747454261913333
$arr[21,5] = 0
$arr[21,6] = 0.07973214163530429
$arr[21,7] = 0
$arr[21,8] = 0.1567843316964894
$arr[21,9] = 0
$arr[21,10] = 0.2916051488777072
$arr[21,11] = 0.386922774581322
$arr[21,12] = 1.647836255007327
$arr[21,13] = 4.55732434792597
$arr[22,0] = 1.632029234738866
$arr[22,1] = 0.1867310688098769
$arr[22,2] = 0.1203619868628607
$arr[22,3] = 0.1224816026024111
$arr[22,4] = 1.746111530706685
$arr[22,5] = 0
$arr[22,6] = 0.07973214163530429
$arr[22,7] = 0
$arr[22,8] = 0.1582472421543408
$arr[22,9] = 0
$arr[22,10] = 0.2853212437653241
$arr[22,11] = 0.3627398154166244
$arr[22,12] = 1.671582260183605
$arr[22,13] = 4.555469510493737
$arr[23,0] = 1.473187919482825
$arr[23,1] = 0.1674387017938557
$arr[23,2] = 0.1197754662594832
$arr[23,3] = 0.1236364558885121
$arr[23,4] = 1.747950084526238
$arr[23,5] = 0
$arr[23,6] = 0.07973214163530429
$arr[23,7] = 0
$arr[23,8] = 0.1599451260517859
$arr[23,9] = 0
$arr[23,10] = 0.2788792273106324
$arr[23,11] = 0.3370413301113686
$arr[23,12] = 1.699415296726702
$arr[23,13] = 4.562594723743416
$ws.Range("B2:O25").Value = $arr
